$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.609230666666667
$ws.Range("H2").Value = 4.827692
$ws.Range("I2").Value = 0.5482851650894511
$ws.Range("J2").Value = 0.5482851650894512
$ws.Range("M2").Value = 70.46291600000001
$ws.Range("N2").Value = 211.388748
$ws.Range("O2").Value = 0.5276750397950939
$ws.Range("P2").Value = 0.5276750397950939
$ws.Range("Q2").Value = 113.3910852899573
$ws.Range("R2").Value = 1020.519767609616
$ws.Range("S2").Value = 0.2893163963076357
$ws.Range("T2").Value = 0.2893163963076358

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.609230666666667
$ws.Range("H3").Value = 4.827692
$ws.Range("I3").Value = 0.5482851650894511
$ws.Range("J3").Value = 0.5482851650894512
$ws.Range("O3").Value = 0.07361176802536967
$ws.Range("P3").Value = 0.07361176802536967
$ws.Range("Q3").Value = 15.81829276925911
$ws.Range("R3").Value = 142.364634923332
$ws.Range("S3").Value = 0.04036024038431619
$ws.Range("T3").Value = 0.0403602403843162

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.609230666666667
$ws.Range("H4").Value = 4.827692
$ws.Range("I4").Value = 0.5482851650894511
$ws.Range("J4").Value = 0.5482851650894512
$ws.Range("M4").Value = 42.505498
$ws.Range("N4").Value = 127.516494
$ws.Range("O4").Value = 0.3183105613832428
$ws.Range("P4").Value = 0.3183105613832428
$ws.Range("Q4").Value = 68.40115088353866
$ws.Range("R4").Value = 615.610357951848
$ws.Range("S4").Value = 0.1745249586977272
$ws.Range("T4").Value = 0.1745249586977272

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.609230666666667
$ws.Range("H5").Value = 4.827692
$ws.Range("I5").Value = 0.5482851650894511
$ws.Range("J5").Value = 0.5482851650894512
$ws.Range("M5").Value = 10.73653933333333
$ws.Range("N5").Value = 32.209618
$ws.Range("O5").Value = 0.08040263079629371
$ws.Range("P5").Value = 0.08040263079629371
$ws.Range("Q5").Value = 17.27756834907289
$ws.Range("R5").Value = 155.498115141656
$ws.Range("S5").Value = 0.04408356969977208
$ws.Range("T5").Value = 0.0440835696997721

$ws.Range("G6").Value = 0.9591736666666666
$ws.Range("I6").Value = 0.3268025542087943
$ws.Range("J6").Value = 0.3268025542087943
$ws.Range("M6").Value = 70.46291600000001
$ws.Range("N6").Value = 211.388748
$ws.Range("O6").Value = 0.5276750397950939
$ws.Range("P6").Value = 0.5276750397950939
$ws.Range("Q6").Value = 67.58617350374533
$ws.Range("R6").Value = 608.275561533708
$ws.Range("S6").Value = 0.1724455507972638
$ws.Range("T6").Value = 0.1724455507972639

$ws.Range("G7").Value = 0.9591736666666666
$ws.Range("I7").Value = 0.3268025542087943
$ws.Range("J7").Value = 0.3268025542087943
$ws.Range("O7").Value = 0.07361176802536967
$ws.Range("P7").Value = 0.07361176802536967
$ws.Range("Q7").Value = 9.428412091676776
$ws.Range("R7").Value = 84.85570882509099
$ws.Range("S7").Value = 0.02405651381051606
$ws.Range("T7").Value = 0.02405651381051607

$ws.Range("G8").Value = 0.9591736666666666
$ws.Range("I8").Value = 0.3268025542087943
$ws.Range("J8").Value = 0.3268025542087943
$ws.Range("M8").Value = 42.505498
$ws.Range("N8").Value = 127.516494
$ws.Range("O8").Value = 0.3183105613832428
$ws.Range("P8").Value = 0.3183105613832428
$ws.Range("Q8").Value = 40.77015437015266
$ws.Range("R8").Value = 366.931389331374
$ws.Range("S8").Value = 0.104024704491679
$ws.Range("T8").Value = 0.104024704491679

$ws.Range("G9").Value = 0.9591736666666666
$ws.Range("I9").Value = 0.3268025542087943
$ws.Range("J9").Value = 0.3268025542087943
$ws.Range("M9").Value = 10.73653933333333
$ws.Range("N9").Value = 32.209618
$ws.Range("O9").Value = 0.08040263079629371
$ws.Range("P9").Value = 0.08040263079629371
$ws.Range("Q9").Value = 10.29820579966422
$ws.Range("R9").Value = 92.68385219697799
$ws.Range("S9").Value = 0.02627578510933545
$ws.Range("T9").Value = 0.02627578510933545

$ws.Range("I10").Value = 0.04674417878325851
$ws.Range("J10").Value = 0.04674417878325852
$ws.Range("M10").Value = 70.46291600000001
$ws.Range("N10").Value = 211.388748
$ws.Range("O10").Value = 0.5276750397950939
$ws.Range("P10").Value = 0.5276750397950939
$ws.Range("Q10").Value = 9.667183248258668
$ws.Range("R10").Value = 87.00464923432801
$ws.Range("S10").Value = 0.02466573639964492
$ws.Range("T10").Value = 0.02466573639964492

$ws.Range("I11").Value = 0.04674417878325851
$ws.Range("J11").Value = 0.04674417878325852
$ws.Range("O11").Value = 0.07361176802536967
$ws.Range("P11").Value = 0.07361176802536967
$ws.Range("S11").Value = 0.003440921645129632
$ws.Range("T11").Value = 0.003440921645129633

$ws.Range("I12").Value = 0.04674417878325851
$ws.Range("J12").Value = 0.04674417878325852
$ws.Range("M12").Value = 42.505498
$ws.Range("N12").Value = 127.516494
$ws.Range("O12").Value = 0.3183105613832428
$ws.Range("P12").Value = 0.3183105613832428
$ws.Range("Q12").Value = 5.831555966609333
$ws.Range("R12").Value = 52.484003699484
$ws.Range("S12").Value = 0.01487916578989769
$ws.Range("T12").Value = 0.01487916578989769

$ws.Range("I13").Value = 0.04674417878325851
$ws.Range("J13").Value = 0.04674417878325852
$ws.Range("M13").Value = 10.73653933333333
$ws.Range("N13").Value = 32.209618
$ws.Range("O13").Value = 0.08040263079629371
$ws.Range("P13").Value = 0.08040263079629371
$ws.Range("Q13").Value = 1.473003092683111
$ws.Range("R13").Value = 13.257027834148
$ws.Range("S13").Value = 0.00375835494858628
$ws.Range("T13").Value = 0.00375835494858628

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.2294253333333333
$ws.Range("H14").Value = 0.688276
$ws.Range("I14").Value = 0.07816810191849585
$ws.Range("J14").Value = 0.07816810191849587
$ws.Range("M14").Value = 70.46291600000001
$ws.Range("N14").Value = 211.388748
$ws.Range("O14").Value = 0.5276750397950939
$ws.Range("P14").Value = 0.5276750397950939
$ws.Range("Q14").Value = 16.16597799093867
$ws.Range("R14").Value = 145.493801918448
$ws.Range("S14").Value = 0.04124735629054926
$ws.Range("T14").Value = 0.04124735629054926

$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.2294253333333333
$ws.Range("H15").Value = 0.688276
$ws.Range("I15").Value = 0.07816810191849585
$ws.Range("J15").Value = 0.07816810191849587
$ws.Range("O15").Value = 0.07361176802536967
$ws.Range("P15").Value = 0.07361176802536967
$ws.Range("Q15").Value = 2.255187628799555
$ws.Range("R15").Value = 20.296688659196
$ws.Range("S15").Value = 0.00575409218540777
$ws.Range("T15").Value = 0.005754092185407772

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.2294253333333333
$ws.Range("H16").Value = 0.688276
$ws.Range("I16").Value = 0.07816810191849585
$ws.Range("J16").Value = 0.07816810191849587
$ws.Range("M16").Value = 42.505498
$ws.Range("N16").Value = 127.516494
$ws.Range("O16").Value = 0.3183105613832428
$ws.Range("P16").Value = 0.3183105613832428
$ws.Range("Q16").Value = 9.751838047149333
$ws.Range("R16").Value = 87.766542424344
$ws.Range("S16").Value = 0.02488173240393896
$ws.Range("T16").Value = 0.02488173240393896

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.2294253333333333
$ws.Range("H17").Value = 0.688276
$ws.Range("I17").Value = 0.07816810191849585
$ws.Range("J17").Value = 0.07816810191849587
$ws.Range("M17").Value = 10.73653933333333
$ws.Range("N17").Value = 32.209618
$ws.Range("O17").Value = 0.08040263079629371
$ws.Range("P17").Value = 0.08040263079629371
$ws.Range("Q17").Value = 2.463234115396444
$ws.Range("R17").Value = 22.169107038568
$ws.Range("S17").Value = 0.00628492103859988
$ws.Range("T17").Value = 0.006284921038599882
